# Apply updated odds values to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("N2").Value = 9

# Row 3
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3

# Row 5
$ws.Range("G5").Value = 2.7
$ws.Range("L5").Value = 3.6

# Row 6
$ws.Range("G6").Value = 1.72
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 4.65
$ws.Range("J6").Value = 2.27
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 4.9
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 2.72
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.65
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.52
$ws.Range("U6").Value = 1.9
$ws.Range("V6").Value = 1.7
$ws.Range("W6").Value = 6
$ws.Range("X6").Value = 7.5
$ws.Range("Y6").Value = 8.25
$ws.Range("Z6").Value = 13.5
$ws.Range("AB6").Value = 32
$ws.Range("AC6").Value = 8.5
$ws.Range("AD6").Value = 6.6
$ws.Range("AE6").Value = 17
$ws.Range("AF6").Value = 90
$ws.Range("AH6").Value = 11.25
$ws.Range("AI6").Value = 26
$ws.Range("AJ6").Value = 15.5
$ws.Range("AK6").Value = 80
$ws.Range("AL6").Value = 50
$ws.Range("AM6").Value = 60
$ws.Range("AN6").Value = 3.5
$ws.Range("AO6").Value = 8.5
$ws.Range("AP6").Value = 18
$ws.Range("AQ6").Value = 29
$ws.Range("AS6").Value = 250
$ws.Range("AT6").Value = 2.5
$ws.Range("AU6").Value = 7.5
$ws.Range("AV6").Value = 70
$ws.Range("AW6").Value = 6.2
$ws.Range("AX6").Value = 27
$ws.Range("AY6").Value = 32
$ws.Range("AZ6").Value = 175
$ws.Range("BA6").Value = 200
$ws.Range("BB6").Value = 450

# Row 7
$ws.Range("N7").Value = 7.8
$ws.Range("S7").Value = 1.36
$ws.Range("T7").Value = 2.9
